$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "BANCAYAN FIESTA DILVER HUMBERTO",
    "LLENQUE ANTON HELEN JOHANA",
    "ANTON INGA FATIMA DEL ROSARIO",
    "FABIANA REBECA ARRUNATEGUI SILUPU",
    "TEMOCHE ECHE URSULA YESSENIA",
    "GONZALES FIESTAS MARIA MARIBEL",
    "VELASCO PEÑA KAREN ARELLYS",
    "HERNANDEZ CARNERO ARTURO SEBASTIAN",
    "FLORES SILUPU MARY CARMEN",
    "RUIZ CHIROQUE CLAUDIA JUDITH",
    "MONDRAGON NONAJULCA MARISOL",
    "ORDINOLA JIBAJA JOSE ALBERTO",
    "BAUTISTA CHAVESTA ERICKA MEDALIT",
    "PINTADO CHASQUERO ESTEFANY",
    "CASTRO ESTRADA CINTHIA PATRICIA",
    "MORENO YANAYACO NAYLA GUADALUPE",
    "PINTADO BENITES CRISTOBAL RODRIGO"
)

$totals = @(75, 70, 66, 65, 63, 62, 55, 52, 49, 49, 48, 48, 46, 45, 34, 32, 1)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}
